# Updates the cryptos list (Coin / Link / Price / Volume(1h)) to the latest
# scraped snapshot, as produced by the scheduled GitHub Actions job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sub3 = [char]0x2083

$rows = @(
    @{Row=2; D='60.411.93'; E='  -2.39%  '},
    @{Row=3; D='2.345.06'; E='  -5.61%  '},
    @{Row=4; E='  -0.07%  '},
    @{Row=5; D='542.97'; E='  -2.41%  '},
    @{Row=6; D='137.50'; E='  -7.13%  '},
    @{Row=7; E='  -0.06%  '},
    @{Row=8; D='0.517'; E='  -14.17%  '},
    @{Row=9; D='2.342.26'; E='  -5.65%  '},
    @{Row=10; E='  -5.22%  '},
    @{Row=11; E='  +0.01%  '},
    @{Row=12; D='5.21'; E='  -5.34%  '},
    @{Row=13; D='0.337'; E='  -6.22%  '},
    @{Row=14; E='  -7.44%  '},
    @{Row=15; D='2.763.87'; E='  -5.76%  '},
    @{Row=16; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.0000161'; E='  -5.21%  '},
    @{Row=17; B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='60.290.61'; E='  -2.41%  '},
    @{Row=18; D='2.342.72'; E='  -5.74%  '},
    @{Row=19; D='10.50'; E='  -6.65%  '},
    @{Row=20; D='4.06'; E='  -4.50%  '},
    @{Row=21; D='313.79'; E='  -3.02%  '},
    @{Row=22; D='6.50'; E='  -10.55%  '},
    @{Row=23; D='0.999'; E='  -0.20%  '},
    @{Row=24; D='1.86'; E='  -1.90%  '},
    @{Row=25; D='62.86'; E='  -2.72%  '},
    @{Row=26; D='8.06'; E='  +2.87%  '},
    @{Row=27; E='  +0.16%  '},
    @{Row=28; D='2.453.60'; E='  -5.95%  '},
    @{Row=29; B='Bittensor'; C='https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; D='507.46'; E='  -11.52%  '},
    @{Row=30; B='PEPE'; C='https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; D=[string]::Concat('0.0', $sub3, '0890'); E='  -12.26%  '},
    @{Row=31; E='  -6.39%  '},
    @{Row=32; D='1.38'; E='  -9.30%  '},
    @{Row=33; E='  -4.99%  '},
    @{Row=34; D='1.80'; E='  -7.23%  '},
    @{Row=35; E='  -6.10%  '},
    @{Row=36; E='  -0.03%  '},
    @{Row=37; D='4.52'; E='  -9.30%  '},
    @{Row=38; D='0.370'; E='  -4.00%  '},
    @{Row=39; B='RenderToken'; C='https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'; D='5.23'; E='  -13.09%  '},
    @{Row=40; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='18.12'; E='  -2.82%  '},
    @{Row=41; D='1.78'; E='  -0.62%  '},
    @{Row=42; D='138.83'; E='  -3.82%  '},
    @{Row=43; E='  -0.06%  '},
    @{Row=44; D='40.01'; E='  -1.66%  '},
    @{Row=45; D='3.53'; E='  -3.61%  '},
    @{Row=46; D='136.95'; E='  -8.31%  '},
    @{Row=47; D='2.06'; E='  -16.51%  '},
    @{Row=48; D='0.0511'; E='  -6.45%  '},
    @{Row=49; D='19.48'; E='  -12.35%  '},
    @{Row=50; D='0.565'; E='  -6.12%  '},
    @{Row=51; D='0.0894'; E='  -5.51%  '}
)

foreach ($r in $rows) {
    $rownum = $r.Row

    if ($r.ContainsKey('B')) {
        $ws.Range('B' + $rownum).Value = $r.B
    }
    if ($r.ContainsKey('C')) {
        $ws.Range('C' + $rownum).Value = $r.C
    }
    if ($r.ContainsKey('D')) {
        # Force the Price column to be stored as text (matching the source
        # data, which is always a text snapshot of the price string) so that
        # values such as "0.370" or "137.50" keep their trailing zeros
        # instead of being reinterpreted as numbers.
        $cellD = $ws.Range('D' + $rownum)
        $cellD.NumberFormat = '@'
        $cellD.Value = $r.D
    }
    if ($r.ContainsKey('E')) {
        $ws.Range('E' + $rownum).Value = $r.E
    }
}
